# Append " (Changed main)" to the end of the first paragraph's text as
# three additional, distinct runs (matching how the document was
# actually authored): " (", "Changed main", ")".
#
# A plain Range.InsertAfter() call merges the newly typed text into the
# neighbouring run whenever the two runs end up with identical
# (absent) run formatting, which would collapse everything back into a
# single <w:r>. Temporarily anchoring a bookmark at the insertion point
# keeps the freshly inserted text in its own run; once the bookmark is
# removed again there is no left-over trace of it in the XML, but the
# run boundary it forced remains.

$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1)
$r = $p1.Range
# Paragraph.Range includes the trailing paragraph mark - back off one
# character so the insertion point sits right after "document." and
# before the pilcrow.
$r.End = $r.End - 1
$r.Collapse(0)

$r.InsertAfter(" (")
$d.Bookmarks.Add("IronEditBoundary1", $r)
$r.Collapse(0)

$r.InsertAfter("Changed main")
$d.Bookmarks.Add("IronEditBoundary2", $r)
$r.Collapse(0)

$r.InsertAfter(")")

$d.Bookmarks("IronEditBoundary1").Delete()
$d.Bookmarks("IronEditBoundary2").Delete()
